$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.90779893785448
$ws.Range("C2").Value = 4.353067274805006
$ws.Range("D2").Value = 8.985354263626654
$ws.Range("F2").Value = 38.73766724238035
$ws.Range("G2").Value = 43.59806339162412
$ws.Range("H2").Value = 17.97269535343644
$ws.Range("J2").Value = 11.23826257278955
$ws.Range("K2").Value = 11.35772170166656
$ws.Range("M2").Value = 16.97011953645116
$ws.Range("N2").Value = 21.19883263867738

$ws.Range("B3").Value = 11.69789897467016
$ws.Range("C3").Value = 4.191309861264235
$ws.Range("D3").Value = 8.957951698739503
$ws.Range("F3").Value = 38.74910436825096
$ws.Range("G3").Value = 43.58077420702121
$ws.Range("H3").Value = 18.01103046508671
$ws.Range("J3").Value = 11.25479606892346
$ws.Range("K3").Value = 11.22480203160382
$ws.Range("M3").Value = 16.9158943202418
$ws.Range("N3").Value = 21.26097190208221

$ws.Range("B4").Value = 11.57011138426196
$ws.Range("C4").Value = 4.089714012164198
$ws.Range("D4").Value = 8.942617520785248
$ws.Range("F4").Value = 38.76455298019199
$ws.Range("G4").Value = 43.58137034467603
$ws.Range("H4").Value = 18.03741704379481
$ws.Range("J4").Value = 11.2664844792197
$ws.Range("K4").Value = 11.14499320506851
$ws.Range("M4").Value = 16.88576291646637
$ws.Range("N4").Value = 21.30099597035992

$ws.Range("B5").Value = 11.51838564156701
$ws.Range("C5").Value = 4.047809911074507
$ws.Range("D5").Value = 8.936748219830289
$ws.Range("F5").Value = 38.77296645710771
$ws.Range("G5").Value = 43.58443187390767
$ws.Range("H5").Value = 18.04888549408609
$ws.Range("J5").Value = 11.27163414047124
$ws.Range("K5").Value = 11.11296170040932
$ws.Range("M5").Value = 16.87428862025749
$ws.Range("N5").Value = 21.31777751770741

$ws.Range("B6").Value = 11.50981985698996
$ws.Range("C6").Value = 4.040823548798962
$ws.Range("D6").Value = 8.935796673727889
$ws.Range("F6").Value = 38.77449139014482
$ws.Range("G6").Value = 43.58511041570993
$ws.Range("H6").Value = 18.05083302926137
$ws.Range("J6").Value = 11.27251258604906
$ws.Range("K6").Value = 11.10767364096009
$ws.Range("M6").Value = 16.87243216164091
$ws.Range("N6").Value = 21.32059258589293

$ws.Range("B7").Value = 11.56941228257971
$ws.Range("C7").Value = 4.089150818660713
$ws.Range("D7").Value = 8.942536823005469
$ws.Range("F7").Value = 38.76465787325822
$ws.Range("G7").Value = 43.58140022359544
$ws.Range("H7").Value = 18.03756881434231
$ws.Range("J7").Value = 11.26655236416385
$ws.Range("K7").Value = 11.14455917912288
$ws.Range("M7").Value = 16.88560490130946
$ws.Range("N7").Value = 21.30122038197929

$ws.Range("B8").Value = 11.8352432021487
$ws.Range("C8").Value = 4.29780528241019
$ws.Range("D8").Value = 8.975599214451083
$ws.Range("F8").Value = 38.73986172936816
$ws.Range("G8").Value = 43.58977524391318
$ws.Range("H8").Value = 17.98532170567198
$ws.Range("J8").Value = 11.24364443764387
$ws.Range("K8").Value = 11.31153927093634
$ws.Range("M8").Value = 16.9507718885904
$ws.Range("N8").Value = 21.21987074468878

$ws.Range("B9").Value = 12.3619727946164
$ws.Range("C9").Value = 4.686206579486374
$ws.Range("D9").Value = 9.052050810098814
$ws.Range("F9").Value = 38.75809087046888
$ws.Range("G9").Value = 43.69510703866259
$ws.Range("H9").Value = 17.90549439079989
$ws.Range("J9").Value = 11.21091035875508
$ws.Range("K9").Value = 11.65161773435791
$ws.Range("M9").Value = 17.10322627125694
$ws.Range("N9").Value = 21.07513285270793

$ws.Range("B10").Value = 12.7480453838794
$ws.Range("C10").Value = 4.955732683293357
$ws.Range("D10").Value = 9.114992387456013
$ws.Range("F10").Value = 38.81217822582329
$ws.Range("G10").Value = 43.8265039372701
$ws.Range("H10").Value = 17.86067448491039
$ws.Range("J10").Value = 11.19428256419442
$ws.Range("K10").Value = 11.90688289564317
$ws.Range("M10").Value = 17.22963587029978
$ws.Range("N10").Value = 20.97773593555449

$ws.Range("B11").Value = 12.92262543080907
$ws.Range("C11").Value = 5.07436173128782
$ws.Range("D11").Value = 9.145024207930039
$ws.Range("F11").Value = 38.84558488692189
$ws.Range("G11").Value = 43.89792325161246
$ws.Range("H11").Value = 17.84329368753885
$ws.Range("J11").Value = 11.18832710347411
$ws.Range("K11").Value = 12.02369589442116
$ws.Range("M11").Value = 17.29011963289073
$ws.Range("N11").Value = 20.93535367176468

$ws.Range("B12").Value = 12.98851624741206
$ws.Range("C12").Value = 5.118673017765928
$ws.Range("D12").Value = 9.156591091459571
$ws.Range("F12").Value = 38.85949582159418
$ws.Range("G12").Value = 43.92663185311091
$ws.Range("H12").Value = 17.83714494848768
$ws.Range("J12").Value = 11.18630291888739
$ws.Range("K12").Value = 12.06798896377434
$ws.Range("M12").Value = 17.31343790006265
$ws.Range("N12").Value = 20.91958017266744

$ws.Range("B13").Value = 12.9743362697191
$ws.Range("C13").Value = 5.109157524808613
$ws.Range("D13").Value = 9.154091416658645
$ws.Range("F13").Value = 38.85644388377963
$ws.Range("G13").Value = 43.92037514917943
$ws.Range("H13").Value = 17.83844992578443
$ws.Range("J13").Value = 11.1867285939072
$ws.Range("K13").Value = 12.05844770035232
$ws.Range("M13").Value = 17.30839769482577
$ws.Range("N13").Value = 20.92296502884629

$ws.Range("B14").Value = 12.92805107455433
$ws.Range("C14").Value = 5.078019695794017
$ws.Range("D14").Value = 9.145971964235395
$ws.Range("F14").Value = 38.84670411812812
$ws.Range("G14").Value = 43.90025185426409
$ws.Range("H14").Value = 17.84277914701079
$ws.Range("J14").Value = 11.18815594451016
$ws.Range("K14").Value = 12.02733897103704
$ws.Range("M14").Value = 17.29202981872359
$ws.Range("N14").Value = 20.93405045547861

$ws.Range("B15").Value = 12.89966957890413
$ws.Range("C15").Value = 5.05886624926328
$ws.Range("D15").Value = 9.141023688004807
$ws.Range("F15").Value = 38.84090221922133
$ws.Range("G15").Value = 43.88814204787256
$ws.Range("H15").Value = 17.84548732017344
$ws.Range("J15").Value = 11.18906031377009
$ws.Range("K15").Value = 12.00829039414026
$ws.Range("M15").Value = 17.28205755850999
$ws.Range("N15").Value = 20.94087648192311

$ws.Range("B16").Value = 12.73660919349406
$ws.Range("C16").Value = 4.947896491323374
$ws.Range("D16").Value = 9.113057326307176
$ws.Range("F16").Value = 38.81017172627408
$ws.Range("G16").Value = 43.82206997255083
$ws.Range("H16").Value = 17.86187092937706
$ws.Range("J16").Value = 11.19470411270018
$ws.Range("K16").Value = 11.89925929478918
$ws.Range("M16").Value = 17.22574188037614
$ws.Range("N16").Value = 20.98054435320801

$ws.Range("B17").Value = 12.63625928227547
$ws.Range("C17").Value = 4.878771687010283
$ws.Range("D17").Value = 9.096254574854621
$ws.Range("F17").Value = 38.79357066939191
$ws.Range("G17").Value = 43.78451258320656
$ws.Range("H17").Value = 17.8726925055564
$ws.Range("J17").Value = 11.19857820911404
$ws.Range("K17").Value = 11.83252006534356
$ws.Range("M17").Value = 17.19194733569627
$ws.Range("N17").Value = 21.00537143967357

$ws.Range("B18").Value = 12.5784460844425
$ws.Range("C18").Value = 4.838640309884638
$ws.Range("D18").Value = 9.086722215012177
$ws.Range("F18").Value = 38.78485105338267
$ws.Range("G18").Value = 43.76400756024834
$ws.Range("H18").Value = 17.87919989557611
$ws.Range("J18").Value = 11.20095791762603
$ws.Range("K18").Value = 11.79420093818686
$ws.Range("M18").Value = 17.17279092602348
$ws.Range("N18").Value = 21.01983248175818

$ws.Range("B19").Value = 12.55885742004119
$ws.Range("C19").Value = 4.824989779817998
$ws.Range("D19").Value = 9.083517617565391
$ws.Range("F19").Value = 38.78204123960829
$ws.Range("G19").Value = 43.75725362337363
$ws.Range("H19").Value = 17.8814517949005
$ws.Range("J19").Value = 11.20178966257948
$ws.Range("K19").Value = 11.78123959285378
$ws.Range("M19").Value = 17.16635362939803
$ws.Range("N19").Value = 21.02475988512269

$ws.Range("B20").Value = 12.64695199618255
$ws.Range("C20").Value = 4.886169001221933
$ws.Range("D20").Value = 9.09802962858571
$ws.Range("F20").Value = 38.79525212874741
$ws.Range("G20").Value = 43.78839716986885
$ws.Range("H20").Value = 17.87151122643709
$ws.Range("J20").Value = 11.19815013440364
$ws.Range("K20").Value = 11.83961788150407
$ws.Range("M20").Value = 17.19551580704807
$ws.Range("N20").Value = 21.00270981151235

$ws.Range("B21").Value = 12.94165262727091
$ws.Range("C21").Value = 5.087182495673916
$ws.Range("D21").Value = 9.148351621072399
$ws.Range("F21").Value = 38.84953075910742
$ws.Range("G21").Value = 43.90611750125771
$ws.Range("H21").Value = 17.84149579536826
$ws.Range("J21").Value = 11.187730429873
$ws.Range("K21").Value = 12.03647509345214
$ws.Range("M21").Value = 17.29682632968625
$ws.Range("N21").Value = 20.930786920718

$ws.Range("B22").Value = 13.13295144384674
$ws.Range("C22").Value = 5.214980791044292
$ws.Range("D22").Value = 9.182370242328554
$ws.Range("F22").Value = 38.89234926678717
$ws.Range("G22").Value = 43.99274498430439
$ws.Range("H22").Value = 17.8244029335123
$ws.Range("J22").Value = 11.18226693153582
$ws.Range("K22").Value = 12.16545531415008
$ws.Range("M22").Value = 17.36544698093206
$ws.Range("N22").Value = 20.88538798036568

$ws.Range("B23").Value = 13.03099283567845
$ws.Range("C23").Value = 5.147111113012769
$ws.Range("D23").Value = 9.164112739989589
$ws.Range("F23").Value = 38.86882622446985
$ws.Range("G23").Value = 43.94562779142944
$ws.Range("H23").Value = 17.83329464679242
$ws.Range("J23").Value = 11.18505982126428
$ws.Range("K23").Value = 12.09659995519034
$ws.Range("M23").Value = 17.32860728443724
$ws.Range("N23").Value = 20.90947152963578

$ws.Range("B24").Value = 12.64211819091185
$ws.Range("C24").Value = 4.882825889134891
$ws.Range("D24").Value = 9.097226728318656
$ws.Range("F24").Value = 38.79448937210628
$ws.Range("G24").Value = 43.78663756072576
$ws.Range("H24").Value = 17.87204439209436
$ws.Range("J24").Value = 11.19834319212312
$ws.Range("K24").Value = 11.83640879948017
$ws.Range("M24").Value = 17.19390165119742
$ws.Range("N24").Value = 21.00391254907019

$ws.Range("B25").Value = 12.21934042520073
$ws.Range("C25").Value = 4.583708991178358
$ws.Range("D25").Value = 9.030155592565457
$ws.Range("F25").Value = 38.74600420609866
$ws.Range("G25").Value = 43.65710512554944
$ws.Range("H25").Value = 17.92466342641598
$ws.Range("J25").Value = 11.21846162918806
$ws.Range("K25").Value = 11.5585010029925
$ws.Range("M25").Value = 17.05940540505073
$ws.Range("N25").Value = 21.11271272834911
